$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').Value = '69.941.23'
$ws.Range('E2').Value = '  +0.70%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').Value = '3.768.49'
$ws.Range('E3').Value = '  -0.02%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.37%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.85'
$ws.Range('E5').Value = '  +3.98%  '

# Row 6: 'Solana' -> 'Solana'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '187.68'
$ws.Range('E6').Value = '  +15.90%  '

# Row 7: 'LidoStakedEther' -> 'LidoStakedEther'
$ws.Range('D7').Value = '3.759.92'
$ws.Range('E7').Value = '  -0.18%  '

# Row 8: 'XRP' -> 'XRP'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.639'
$ws.Range('E8').Value = '  -1.88%  '

# Row 9: 'USDC' -> 'USDC'
$ws.Range('E9').Value = '  -0.51%  '

# Row 10: 'Cardano' -> 'Cardano'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.733'
$ws.Range('E10').Value = '  +0.89%  '

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.165'
$ws.Range('E11').Value = '  -2.76%  '

# Row 12: 'Avalanche' -> 'Avalanche'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '57.99'
$ws.Range('E12').Value = '  +12.63%  '

# Row 13: 'ShibaInu' -> 'ShibaInu'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000299'
$ws.Range('E13').Value = '  -4.46%  '

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.98'
$ws.Range('E14').Value = '  +0.80%  '

# Row 15: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range('D15').Value = '4.361.94'
$ws.Range('E15').Value = '  -0.40%  '

# Row 16: 'WrappedEther' -> 'WrappedEther'
$ws.Range('D16').Value = '3.777.51'
$ws.Range('E16').Value = '  -0.68%  '

# Row 17: 'Chainlink' -> 'Chainlink'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.71'
$ws.Range('E17').Value = '  -2.93%  '

# Row 18: 'Uniswap' -> 'Uniswap'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.07'
$ws.Range('E18').Value = '  -2.50%  '

# Row 19: 'Polygon' -> 'TRON'
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.127'
$ws.Range('E19').Value = '  -1.37%  '

# Row 20: 'TRON' -> 'Polygon'
$ws.Range('B20').Value = 'Polygon'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.14'
$ws.Range('E20').Value = '  -3.19%  '

# Row 21: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range('D21').Value = '69.666.93'
$ws.Range('E21').Value = '  +0.33%  '

# Row 22: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '417.41'
$ws.Range('E22').Value = '  -2.68%  '

# Row 23: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.68'
$ws.Range('E23').Value = '  +1.68%  '

# Row 24: 'Litecoin' -> 'Litecoin'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '90.29'
$ws.Range('E24').Value = '  -1.38%  '

# Row 25: 'ImmutableX' -> 'ImmutableX'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.09'
$ws.Range('E25').Value = '  -3.04%  '

# Row 26: 'RenderToken' -> 'InternetComputer(DFINITY)'
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.09'
$ws.Range('E26').Value = '  -3.27%  '

# Row 27: 'InternetComputer(DFINITY)' -> 'RenderToken'
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.27'
$ws.Range('E27').Value = '  +2.43%  '

# Row 28: 'Toncoin' -> 'Toncoin'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.02'
$ws.Range('E28').Value = '  +4.12%  '

# Row 29: 'LEO' -> 'LEO'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.09'
$ws.Range('E29').Value = '  +2.64%  '

# Row 30: 'Filecoin' -> 'Filecoin'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.65'
$ws.Range('E30').Value = '  -4.76%  '

# Row 31: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.38'
$ws.Range('E31').Value = '  -2.26%  '

# Row 32: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.49'
$ws.Range('E32').Value = '  -4.24%  '

# Row 33: 'Cosmos' -> 'Cosmos'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.65'
$ws.Range('E33').Value = '  -3.74%  '

# Row 34: 'Hedera' -> 'Hedera'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.120'
$ws.Range('E34').Value = '  -1.29%  '

# Row 35: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '44.74'
$ws.Range('E35').Value = '  -4.69%  '

# Row 36: 'Bittensor' -> 'OKB'
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '65.37'
$ws.Range('E36').Value = '  -3.60%  '

# Row 37: 'OKB' -> 'Bittensor'
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '610.50'
$ws.Range('E37').Value = '  -1.38%  '

# Row 38: 'PEPE' -> 'PEPE'
$ws.Range('D38').Value = '0.0₃0919'
$ws.Range('E38').Value = '  -3.85%  '

# Row 39: 'TheGraph' -> 'TheGraph'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.412'
$ws.Range('E39').Value = '  -0.48%  '

# Row 40: 'Dai' -> 'Dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.11%  '

# Row 41: 'FirstDigitalUSD' -> 'FirstDigitalUSD'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.44%  '

# Row 42: 'Kaspa' -> 'Kaspa'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.139'
$ws.Range('E42').Value = '  -1.40%  '

# Row 43: 'ThetaToken' -> 'ThetaToken'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.10'
$ws.Range('E43').Value = '  -1.72%  '

# Row 44: 'dogwifhat' -> 'Fetch.AI'
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  +2.62%  '

# Row 45: 'Fetch.AI' -> 'dogwifhat'
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.08'
$ws.Range('E45').Value = '  +0.49%  '

# Row 46: 'VeChain' -> 'VeChain'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0450'
$ws.Range('E46').Value = '  -1.56%  '

# Row 47: 'THORChain' -> 'THORChain'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.46'
$ws.Range('E47').Value = '  -2.15%  '

# Row 48: 'WEMIXToken' -> 'ApeXProtocol'
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.26'
$ws.Range('E48').Value = '  +1.30%  '

# Row 49: 'ApeXProtocol' -> 'Stellar'
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.136'
$ws.Range('E49').Value = '  -2.81%  '

# Row 50: 'Stellar' -> 'WEMIXToken'
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.75'
$ws.Range('E50').Value = '  -1.73%  '

# Row 51: 'Maker' -> 'Maker'
$ws.Range('D51').Value = '2.807.32'
$ws.Range('E51').Value = '  +0.11%  '
